# fdo#51601 test fixture update
#  - note the bug number next to the existing sample data on Sheet1
#  - add Sheet2 with a formula exercising AVERAGEIF (the function fdo#51601
#    is about), right after Sheet1
#  - leave the cursor/selection the way the author left them

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Reference the bug right next to the rest of the sample data.
$ws1.Range("B3").Value = "Fdo#51601"

# Re-stamp the row heights on Sheet1 (values unchanged, just re-applied).
$ws1.Rows.Item(1).RowHeight = 12.8
$ws1.Rows.Item(2).RowHeight = 12.8
$ws1.Rows.Item(3).RowHeight = 12.8
$ws1.Rows.Item(4).RowHeight = 12.8
$ws1.Rows.Item(5).RowHeight = 12.8

# Add Sheet2 right after Sheet1 and put the AVERAGEIF repro formula in it.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"
$ws2.Range("C1").Formula = "=averageif(A2:B2)"
$ws2.Rows.Item(1).RowHeight = 12.1
$ws2.Range("C1").Select() | Out-Null

# Leave Sheet1 active, with the selection where the author left it.
$ws1.Activate() | Out-Null
$ws1.Range("D4").Select() | Out-Null
